# The deck currently uses the "Integral" (Red Violet) theme for the
# slides/slide master (ppt/theme/theme2.xml, wired to the slide master +
# presentation) while the notes master keeps the stock "Office Theme"
# (ppt/theme/theme1.xml). The authored change swaps the two theme
# payloads so the slides/slide master fall back to the plain "Office
# Theme" palette. Re-apply the 12 standard Office-theme colors (in
# msoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) to the live theme that backs the slide master/presentation,
# via PowerPoint's theme-color COM surface.

$p = $ppt.ActivePresentation

# Values below are plain VBA-style RGB(r,g,b) longs (r + g*256 + b*65536)
# for the stock "Office Theme" color scheme.
$officeThemeColors = @(
    0x000000,  # 1  dk1      000000 -> r=00 g=00 b=00
    0xFFFFFF,  # 2  lt1      FFFFFF -> r=FF g=FF b=FF
    0x6A5444,  # 3  dk2      44546A -> r=44 g=54 b=6A
    0xE6E6E7,  # 4  lt2      E7E6E6 -> r=E7 g=E6 b=E6
    0xD59B5B,  # 5  accent1  5B9BD5 -> r=5B g=9B b=D5
    0x317DED,  # 6  accent2  ED7D31 -> r=ED g=7D b=31
    0xA5A5A5,  # 7  accent3  A5A5A5 -> r=A5 g=A5 b=A5
    0x00C0FF,  # 8  accent4  FFC000 -> r=FF g=C0 b=00
    0xC47244,  # 9  accent5  4472C4 -> r=44 g=72 b=C4
    0x47AD70,  # 10 accent6  70AD47 -> r=70 g=AD b=47
    0xC16305,  # 11 hlink    0563C1 -> r=05 g=63 b=C1
    0x724F95   # 12 folHlink 954F72 -> r=95 g=4F b=72
)

# The theme is shared by every slide, so touching it through the first
# slide's ThemeColorScheme updates the one live theme part for the
# whole deck (slide master + all slides + all layouts).
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.ThemeColorScheme.Count; $i++) {
    $s.ThemeColorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}

# Keep the slide master's classic 8-slot color scheme (background/text/
# shadow/title/fill/accent1-3) in sync with the same palette.
$m = $p.SlideMaster
for ($i = 1; $i -le $m.ColorScheme.Count; $i++) {
    $m.ColorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
